$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. Insert a new row at
# position 128 (pushing the existing rows 128-144 down to 129-145) and fill
# it with the new observation's values.
$ws.Rows(128).Insert()

$ws.Range("A128").Value2 = 6
$ws.Range("B128").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C128").Value2 = "Metropolitana"
$ws.Range("D128").Value2 = 45154
$ws.Range("E128").Value2 = 13
$ws.Range("F128").Value2 = 100114007
$ws.Range("G128").Value2 = "Jengibre"
$ws.Range("H128").Value2 = "Sin especificar"
$ws.Range("I128").Value2 = "Primera"
$ws.Range("J128").Value2 = 560
$ws.Range("K128").Value2 = 14000
$ws.Range("L128").Value2 = 15000
$ws.Range("M128").Value2 = 14571
$ws.Range("N128").Value2 = "$/caja 13 kilos"
$ws.Range("O128").Value2 = "Perú"
$ws.Range("P128").Value2 = 1121
$ws.Range("Q128").Value2 = 13
$ws.Range("R128").Value2 = "Hortaliza"

# Preserve the date-formatted style for the new row's Fecha cell, matching
# the other rows in the column.
$ws.Range("D128").NumberFormat = $ws.Range("D129").NumberFormat
